# Apply updated cryptocurrency price/volume data (Wed Oct 11 05:05:19 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.171.43"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  -1.79%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.562.08"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -1.61%  "
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "206.72"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -0.12%  "
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -1.43%  "
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "22.35"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +0.66%  "
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -1.82%  "
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +0.19%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0860"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -0.83%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.785.78"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -1.54%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.567.09"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -1.07%  "
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -2.31%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "62.85"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "27.159.38"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  -1.77%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "212.65"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -2.93%  "
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -1.24%  "
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -1.18%  "
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  -0.04%  "
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -0.30%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.36"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -2.42%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.98"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "152.26"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -0.74%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "6.59"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -3.51%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "14.88"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -1.05%  "
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.15"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.62%  "
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -0.87%  "
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -1.57%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.381.25"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  +0.87%  "
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.58%  "
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +0.65%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.945"
$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -3.12%  "
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -1.04%  "
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -1.06%  "
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -0.82%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.518"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  -3.23%  "
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.989"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +1.50%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.80"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  +5.04%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "63.43"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  -0.93%  "
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  +0.02%  "
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  +0.67%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "1.697.72"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  -1.57%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "85.62"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -2.18%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0₇0993"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  -1.16%  "
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  -0.60%  "
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +0.18%  "
